$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.072439420521087
$ws.Range("D2").Value = 1.063407812128526
$ws.Range("E2").Value = 1.085238353133776
$ws.Range("F2").Value = 1.090492257209315
$ws.Range("I2").Value = 1.045276568057868
$ws.Range("J2").Value = 1.077358960597939
$ws.Range("K2").Value = 1.066126533555585
$ws.Range("L2").Value = 1.08789930005159
$ws.Range("M2").Value = 1.093139672317776
$ws.Range("N2").Value = 1.078888933896788

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.074347983584873
$ws.Range("D3").Value = 1.064809701476444
$ws.Range("E3").Value = 1.087119297462459
$ws.Range("F3").Value = 1.092395581855102
$ws.Range("I3").Value = 1.045741837743584
$ws.Range("J3").Value = 1.078921349851871
$ws.Range("K3").Value = 1.067341955390335
$ws.Range("L3").Value = 1.089596848288
$ws.Range("M3").Value = 1.094860546147278
$ws.Range("N3").Value = 1.08045354192267

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.075579358946556
$ws.Range("D4").Value = 1.065713497702489
$ws.Range("E4").Value = 1.088333142776576
$ws.Range("F4").Value = 1.093623915174329
$ws.Range("I4").Value = 1.046039999091029
$ws.Range("J4").Value = 1.0799284569511
$ws.Range("K4").Value = 1.068124571673167
$ws.Range("L4").Value = 1.090691588875891
$ws.Range("M4").Value = 1.095970387931743
$ws.Range("N4").Value = 1.081462079229496

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07609618877062
$ws.Range("D5").Value = 1.066092672348484
$ws.Range("E5").Value = 1.08884268346864
$ws.Range("F5").Value = 1.09413954835116
$ws.Range("I5").Value = 1.046164657605474
$ws.Range("J5").Value = 1.080350936156249
$ws.Range("K5").Value = 1.068452676325648
$ws.Range("L5").Value = 1.091150952518573
$ws.Range("M5").Value = 1.096436102210427
$ws.Range("N5").Value = 1.081885158403575

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.076182918060409
$ws.Range("D6").Value = 1.066156292014724
$ws.Range("E6").Value = 1.088928193568012
$ws.Range("F6").Value = 1.094226081512171
$ws.Range("I6").Value = 1.04618554810432
$ws.Range("J6").Value = 1.08042181946071
$ws.Range("K6").Value = 1.068507713708391
$ws.Range("L6").Value = 1.091228031452779
$ws.Range("M6").Value = 1.096514247572383
$ws.Range("N6").Value = 1.081956142370459

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.075586268127305
$ws.Range("D7").Value = 1.065718567304151
$ws.Range("E7").Value = 1.088339954244831
$ws.Range("F7").Value = 1.093630808043624
$ws.Range("I7").Value = 1.046041667481968
$ws.Range("J7").Value = 1.07993410568705
$ws.Range("K7").Value = 1.068128959370978
$ws.Range("L7").Value = 1.090697730291872
$ws.Range("M7").Value = 1.09597661419676
$ws.Range("N7").Value = 1.0814677359873

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.073085184116401
$ws.Range("D8").Value = 1.063882282782478
$ws.Range("E8").Value = 1.085874711512616
$ws.Range("F8").Value = 1.091136177898669
$ws.Range("I8").Value = 1.04543441180696
$ws.Range("J8").Value = 1.077887787093137
$ws.Range("K8").Value = 1.066538094482336
$ws.Range("L8").Value = 1.088473768224227
$ws.Range("M8").Value = 1.093722021942489
$ws.Range("N8").Value = 1.079418511386267

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.068649551247858
$ws.Range("D9").Value = 1.060620477985679
$ws.Range("E9").Value = 1.081504900619634
$ws.Range("F9").Value = 1.086714623652229
$ws.Range("I9").Value = 1.044341893043063
$ws.Range("J9").Value = 1.074251591221536
$ws.Range("K9").Value = 1.063704778317551
$ws.Range("L9").Value = 1.084525847327366
$ws.Range("M9").Value = 1.08972017960147
$ws.Range("N9").Value = 1.075777151699451

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.065672116885626
$ws.Range("D10").Value = 1.058427576749098
$ws.Range("E10").Value = 1.078573212894531
$ws.Range("F10").Value = 1.083748456259397
$ws.Range("I10").Value = 1.043598093477736
$ws.Range("J10").Value = 1.071806042887199
$ws.Range("K10").Value = 1.061794912566939
$ws.Range("L10").Value = 1.081873294127701
$ws.Range("M10").Value = 1.087031702740644
$ws.Range("N10").Value = 1.07332813040597

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.064377749391772
$ws.Range("D11").Value = 1.057473477639326
$ws.Range("E11").Value = 1.077299113001285
$ws.Range("F11").Value = 1.082459427451707
$ws.Range("I11").Value = 1.043272274753587
$ws.Range("J11").Value = 1.070741782161458
$ws.Range("K11").Value = 1.060962759222648
$ws.Range("L11").Value = 1.08071957584303
$ws.Range("M11").Value = 1.085862432212702
$ws.Range("N11").Value = 1.072262358307926

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.063896170398072
$ws.Range("D12").Value = 1.057118381154918
$ws.Range("E12").Value = 1.076825132667818
$ws.Range("F12").Value = 1.081979901781231
$ws.Range("I12").Value = 1.043150681016501
$ws.Range("J12").Value = 1.070345648683402
$ws.Range("K12").Value = 1.060652868019901
$ws.Range("L12").Value = 1.080290239156177
$ws.Range("M12").Value = 1.085427318948544
$ws.Range("N12").Value = 1.071865662274894

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.063999507016831
$ws.Range("D13").Value = 1.057194582576054
$ws.Range("E13").Value = 1.076926836140897
$ws.Range("F13").Value = 1.08208279475489
$ws.Range("I13").Value = 1.043176789199546
$ws.Range("J13").Value = 1.070430658075354
$ws.Range("K13").Value = 1.060719376849042
$ws.Range("L13").Value = 1.08038236962424
$ws.Range("M13").Value = 1.085520688505782
$ws.Range("N13").Value = 1.071950792389933

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.064337958245846
$ws.Range("D14").Value = 1.057444139663134
$ws.Range("E14").Value = 1.077259948519213
$ws.Range("F14").Value = 1.082419804593158
$ws.Range("I14").Value = 1.043262235449025
$ws.Range("J14").Value = 1.070709054482395
$ws.Range("K14").Value = 1.060937159776229
$ws.Range("L14").Value = 1.080684103067409
$ws.Range("M14").Value = 1.085826481940246
$ws.Range("N14").Value = 1.072229584151805

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.064546383371159
$ws.Range("D15").Value = 1.057597806548083
$ws.Range("E15").Value = 1.077465093569749
$ws.Range("F15").Value = 1.082627350980887
$ws.Range("I15").Value = 1.043314805942037
$ws.Range("J15").Value = 1.070880474479508
$ws.Range("K15").Value = 1.061071237547462
$ws.Range("L15").Value = 1.080869905098553
$ws.Range("M15").Value = 1.086014785488228
$ws.Range("N15").Value = 1.072401247584976

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.065757909885154
$ws.Range("D16").Value = 1.05849079970029
$ws.Range("E16").Value = 1.078657670590867
$ws.Range("F16").Value = 1.083833904740728
$ws.Range("I16").Value = 1.043619637421845
$ws.Range("J16").Value = 1.071876560528213
$ws.Range("K16").Value = 1.061850029594882
$ws.Range("L16").Value = 1.081949752391809
$ws.Range("M16").Value = 1.087109193198847
$ws.Range("N16").Value = 1.073398748190123

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.066516481008357
$ws.Range("D17").Value = 1.059049718547169
$ws.Range("E17").Value = 1.079404479076762
$ws.Range("F17").Value = 1.084589480577689
$ws.Range("I17").Value = 1.043809841657059
$ws.Range("J17").Value = 1.072499940015983
$ws.Range("K17").Value = 1.062337150305319
$ws.Range("L17").Value = 1.0826257195702
$ws.Range("M17").Value = 1.087794294372474
$ws.Range("N17").Value = 1.074023012948269

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.066958450476561
$ws.Range("D18").Value = 1.059375287777098
$ws.Range("E18").Value = 1.079839631788895
$ws.Range("F18").Value = 1.085029747250077
$ws.Range("I18").Value = 1.043920423359342
$ws.Range("J18").Value = 1.072863034581827
$ws.Range("K18").Value = 1.062620781933529
$ws.Range("L18").Value = 1.083019504802257
$ws.Range("M18").Value = 1.088193407608206
$ws.Range("N18").Value = 1.074386623150043

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.067109067880993
$ws.Range("D19").Value = 1.059486224566999
$ws.Range("E19").Value = 1.079987932308814
$ws.Range("F19").Value = 1.085179791524175
$ws.Range("I19").Value = 1.043958067804043
$ws.Range("J19").Value = 1.072986754180422
$ws.Range("K19").Value = 1.062717409056431
$ws.Range("L19").Value = 1.083153692072998
$ws.Range("M19").Value = 1.088329411647909
$ws.Range("N19").Value = 1.07451051844466

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.066435144594636
$ws.Range("D20").Value = 1.058989797339612
$ws.Range("E20").Value = 1.079324400101642
$ws.Range("F20").Value = 1.084508460938739
$ws.Range("I20").Value = 1.043789471954599
$ws.Range("J20").Value = 1.072433110414478
$ws.Range("K20").Value = 1.062284938450606
$ws.Range("L20").Value = 1.08255324604179
$ws.Range("M20").Value = 1.087720840823524
$ws.Range("N20").Value = 1.073956088441063

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.064238314889071
$ws.Range("D21").Value = 1.057370670824057
$ws.Range("E21").Value = 1.077161875366902
$ws.Range("F21").Value = 1.082320583739883
$ws.Range("I21").Value = 1.043237089453795
$ws.Range("J21").Value = 1.070627096468403
$ws.Range("K21").Value = 1.060873050134584
$ws.Range("L21").Value = 1.080595272189468
$ws.Range("M21").Value = 1.085736455491842
$ws.Range("N21").Value = 1.072147509748032

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.062852478980626
$ws.Range("D22").Value = 1.056348594056532
$ws.Range("E22").Value = 1.075798016211461
$ws.Range("F22").Value = 1.080940783945457
$ws.Range("I22").Value = 1.042886482382643
$ws.Range("J22").Value = 1.069486831397473
$ws.Range("K22").Value = 1.059980747337196
$ws.Range("L22").Value = 1.079359609700511
$ws.Range("M22").Value = 1.084484188208067
$ws.Range("N22").Value = 1.071005625369912

$ws.Range("B23").Value = 1.019999999999999
$ws.Range("C23").Value = 1.063587581196538
$ws.Range("D23").Value = 1.056890807506156
$ws.Range("E23").Value = 1.076521428999102
$ws.Range("F23").Value = 1.081672647294332
$ws.Range("I23").Value = 1.043072661193883
$ws.Range("J23").Value = 1.07009176472871
$ws.Range("K23").Value = 1.060454214598125
$ws.Range("L23").Value = 1.080015101614563
$ws.Range("M23").Value = 1.085148482604204
$ws.Range("N23").Value = 1.071611417775866

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.066471898539293
$ws.Range("D24").Value = 1.059016874509191
$ws.Range("E24").Value = 1.079360585727449
$ws.Range("F24").Value = 1.084545071609688
$ws.Range("I24").Value = 1.043798677262318
$ws.Range("J24").Value = 1.07246330941678
$ws.Range("K24").Value = 1.062308532277658
$ws.Range("L24").Value = 1.082585995238684
$ws.Range("M24").Value = 1.087754032848587
$ws.Range("N24").Value = 1.073986330329413

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.069799762247797
$ws.Range("D25").Value = 1.061466905969302
$ws.Range("E25").Value = 1.082637771723748
$ws.Range("F25").Value = 1.08786086798685
$ws.Range("I25").Value = 1.044627032830209
$ws.Range("J25").Value = 1.075195333623935
$ws.Range("K25").Value = 1.064440897630089
$ws.Range("L25").Value = 1.085550030156053
$ws.Range("M25").Value = 1.090758297335426
$ws.Range("N25").Value = 1.076722234324309
